$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2; this pushes the existing "TUM_Info"/"Yes" row
# down to row 3 (it keeps referencing the same shared-string slot).
$ws.Rows("2:2").Insert()

# The program that used to be in row 2 (now row 3) is renamed.
$ws.Range("A3").Value = "RWTH_Aachen_Data Science"
$ws.Range("B3").Value = "Yes"

# Add the new programs under it.
$ws.Range("A4").Value = "Freie Uni Berlin - Data Science"
$ws.Range("B4").Value = "Yes"

$ws.Range("A5").Value = "TU Berlin Computer Science"
$ws.Range("B5").Value = "Yes"

# Fill the freshly-inserted row 2 with the renamed TUM program.
$ws.Range("A2").Value = "TUM_Informatics"
$ws.Range("B2").Value = "Yes"

# Append one more program as a new row 6.
$ws.Range("A6").Value = "TU Data Engineering and Analytics"
$ws.Range("B6").Value = "Yes"

# Extend the Yes/No dropdown validation down to the new last row.
$ws.Range("B1:B6").Validation.Delete()
$ws.Range("B1:B6").Validation.Add(3, 1, 1, '"Yes,No"')

# Match the saved selection state.
[void]$ws.Range("C3").Select()
